# Weekly update of fruit/vegetable price data.
# The data rows (2-16) are re-shuffled with a new weekly dataset.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Categoría ID,
#          G Categoría, H Variedad, I Calidad, J Volumen, K Precio mínimo,
#          L Precio máximo, M Precio promedio ponderado, N Unidad de comercialización,
#          O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificación

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44467, 13, 100112013, "Alcachofa", "Española", "Primera", 35, 12000, 12000, 12000, "`$/caja 30 unidades", "Provincia de Limarí", 400, 30, "Hortaliza")
    3  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44425, 13, 100112013, "Alcachofa", "Española", "Primera", 35, 14000, 14000, 14000, "`$/caja 30 unidades", "Provincia de Limarí", 467, 30, "Hortaliza")
    4  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44449, 13, 100112013, "Alcachofa", "Española", "Primera", 45, 12000, 12000, 12000, "`$/caja 30 unidades", "Provincia de Limarí", 400, 30, "Hortaliza")
    5  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44453, 13, 100112013, "Alcachofa", "Española", "Primera", 50, 12000, 12000, 12000, "`$/caja 30 unidades", "Provincia de Limarí", 400, 30, "Hortaliza")
    6  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44418, 13, 100112013, "Alcachofa", "Española", "Primera", 30, 15000, 15000, 15000, "`$/caja 30 unidades", "Provincia de Limarí", 500, 30, "Hortaliza")
    7  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44474, 13, 100112013, "Alcachofa", "Española", "Primera", 45, 10000, 10000, 10000, "`$/caja 30 unidades", "Provincia de Limarí", 333, 30, "Hortaliza")
    8  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44460, 13, 100112013, "Alcachofa", "Española", "Primera", 45, 13000, 13000, 13000, "`$/caja 30 unidades", "Provincia de Limarí", 433, 30, "Hortaliza")
    9  = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44432, 13, 100112013, "Alcachofa", "Española", "Primera", 25, 14000, 14000, 14000, "`$/caja 30 unidades", "Provincia del Elquí", 467, 30, "Hortaliza")
    10 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44446, 13, 100112013, "Alcachofa", "Española", "Primera", 25, 14000, 14000, 14000, "`$/caja 30 unidades", "Provincia de Limarí", 467, 30, "Hortaliza")
    11 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44841, 13, 100112013, "Alcachofa", "Española", "Primera", 45, 12000, 12000, 12000, "`$/caja 30 unidades", "Provincia de Limarí", 400, 30, "Hortaliza")
    12 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44841, 13, 100112013, "Alcachofa", "Española", "Segunda", 45, 10000, 10000, 10000, "`$/caja 40 unidades", "Provincia de Limarí", 250, 40, "Hortaliza")
    13 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44435, 13, 100112013, "Alcachofa", "Española", "Primera", 25, 14000, 14000, 14000, "`$/caja 30 unidades", "Provincia de Limarí", 467, 30, "Hortaliza")
    14 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44435, 13, 100112013, "Alcachofa", "Española", "Primera", 25, 14000, 14000, 14000, "`$/caja 30 unidades", "Provincia del Elquí", 467, 30, "Hortaliza")
    15 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44421, 13, 100112013, "Alcachofa", "Española", "Primera", 25, 15000, 16000, 15400, "`$/caja 30 unidades", "Provincia de Limarí", 513, 30, "Hortaliza")
    16 = @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44376, 13, 100112013, "Alcachofa", "Española", "Primera", 25, 18000, 18000, 18000, "`$/caja 30 unidades", "Provincia de Limarí", 600, 30, "Hortaliza")
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    for ($colIdx = 0; $colIdx -lt $rowValues.Length; $colIdx++) {
        $col = $colIdx + 1
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$colIdx]
    }
}
